$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark from the first paragraph; it will be
# re-added at the end of the new final paragraph below.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Insert the new paragraphs (per-icon headings + their attribution blocks +
# blank separator paragraphs) right before the paragraph mark that ends the
# first paragraph ("Meditation Icon:"). The very last paragraph in this
# block is a duplicate of the original attribution paragraph with the
# relocated _GoBack bookmark appended to it.
$p1 = $d.Paragraphs(1).Range
$insertionPoint = $d.Range($p1.End - 1, $p1.End - 1)

$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">&lt;div&gt;Icons made by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="https://www.freepik.com/" title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Freepik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>"&gt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Freepik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;/a&gt; from &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="https://www.flaticon.com/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flaticon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">"&gt;www.flaticon.com&lt;/a&gt; is licensed by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="http://creativecommons.org/licenses/by/3.0/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="Creative Commons BY 3.0" target="_blank"&gt;CC 3.0 BY&lt;/a&gt;&lt;/div&gt;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Uniform Icon:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">&lt;div&gt;Icons made by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="https://www.flaticon.com/authors/smashicons" title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Smashicons</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>"&gt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Smashicons</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;/a&gt; from &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="https://www.flaticon.com/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flaticon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">"&gt;www.flaticon.com&lt;/a&gt; is licensed by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="http://creativecommons.org/licenses/by/3.0/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="Creative Commons BY 3.0" target="_blank"&gt;CC 3.0 BY&lt;/a&gt;&lt;/div&gt;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Kicking Icon:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">&lt;div&gt;Icons made by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="https://www.freepik.com/" title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Freepik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>"&gt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Freepik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;/a&gt; from &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="https://www.flaticon.com/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flaticon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">"&gt;www.flaticon.com&lt;/a&gt; is licensed by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="http://creativecommons.org/licenses/by/3.0/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="Creative Commons BY 3.0" target="_blank"&gt;CC 3.0 BY&lt;/a&gt;&lt;/div&gt;</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Yin Yang: </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">&lt;div&gt;Icons made by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="https://www.freepik.com/" title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Freepik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>"&gt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Freepik</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&lt;/a&gt; from &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="https://www.flaticon.com/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Flaticon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">"&gt;www.flaticon.com&lt;/a&gt; is licensed by &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">="http://creativecommons.org/licenses/by/3.0/" </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">    title="Creative Commons BY 3.0" target="_blank"&gt;CC 3.0 BY&lt;/a&gt;&lt;/div&gt;</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$insertionPoint.InsertXML($newParagraphsXml)

# The original attribution paragraph (now pushed to the very end of the
# document) is redundant with the newly inserted "Yin Yang" attribution
# paragraph above, so remove it.
$trailing = $d.Paragraphs($d.Paragraphs.Count)
$trailing.Range.Delete()
